$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 2214149.3
$ws.Range("E2").Value = 2214149.3
$ws.Range("G2").Value = 2214149.3
$ws.Range("I2").Value = 1761181.52
$ws.Range("J2").Value = 239578.26
$ws.Range("L2").Value = 2000759.78
$ws.Range("O2").Value = 213389.48

$ws.Range("I3").Value = 2115213.29
$ws.Range("L3").Value = 1553731.43
$ws.Range("O3").Value = 360696.51

$ws.Range("I4").Value = 1235426.37
$ws.Range("L4").Value = 983809.16
$ws.Range("O4").Value = 761918.39

$ws.Range("I5").Value = 236352.14
$ws.Range("L5").Value = 284157.14
$ws.Range("O5").Value = 50477.83

$ws.Range("I6").Value = 9831338.699999999
$ws.Range("L6").Value = 11999348.1
$ws.Range("O6").Value = 9680745.9

$ws.Range("I7").Value = 54929.84
$ws.Range("L7").Value = 82394.75999999999
$ws.Range("O7").Value = 466903.44

$ws.Range("I8").Value = 95780.64
$ws.Range("L8").Value = 143670.96
$ws.Range("O8").Value = 574684.04

$ws.Range("I9").Value = 8920
$ws.Range("L9").Value = 13380
$ws.Range("O9").Value = 53520

$ws.Range("I10").Value = 60785.5
$ws.Range("L10").Value = 91178.25
$ws.Range("O10").Value = 516676.75

$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 1402841.25
$ws.Range("G11").Value = 1402841.25
$ws.Range("I11").Value = 187045.52
$ws.Range("J11").Value = 93522.75999999999
$ws.Range("L11").Value = 280568.28
$ws.Range("O11").Value = 1122272.97

$ws.Range("I12").Value = 75193.75999999999
$ws.Range("L12").Value = 112790.64
$ws.Range("O12").Value = 263178.11

$ws.Range("I13").Value = 473327.98
$ws.Range("L13").Value = 709991.97
$ws.Range("O13").Value = 4023287.28

$ws.Range("I14").Value = 636143.62
$ws.Range("L14").Value = 954222.89
$ws.Range("O14").Value = 3816966.61
